# Adding cost columns to target model
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1: new header cells T1:AC1 ---
$ws.Range("T1").Value  = "LA.Avg.Gross.Weekly.Cost.Per.Person"
$ws.Range("U1").Value  = "LA.Cost.Of.Care.18-64"
$ws.Range("V1").Value  = "LA.Cost.Of.Care.65-74"
$ws.Range("W1").Value  = "LA.Cost.Of.Care.75-84"
$ws.Range("X1").Value  = "LA.Cost.Of.Care.75pl"
$ws.Range("Y1").Value  = "Other Supplementary LA Cost Metrics"
$ws.Range("Z1").Value  = "…"
$ws.Range("AA1:AC1").Value = ".."

# --- Row 3: sub-header / annotation cells ---
$ws.Range("N3:Q3").Value = "Res / Nursing"
$ws.Range("R3:S3").Value = "Res / Nursing?"
$ws.Range("T3:X3").Value = "Res / Nursing"
$ws.Range("Y3").Value    = "TBD!"

# --- Row 5: free-text notes ---
$ws.Range("A5").Value    = "~150 LAs"
$ws.Range("C5").Value    = "A row for every year between now and 2037"
$ws.Range("F5:M5").Value = "From ONS SRC dataset"
$ws.Range("U5:X5").Value = "Estiamted LA-funded occupancy, x LA Avg Gross Weekly Cost of Care"

# --- Column widths for the newly-introduced columns ---
$ws.Range("U1:X1").ColumnWidth = 20.45
$ws.Range("Y1").ColumnWidth    = 32.45

# --- View state: scroll / selection to match author's final position ---
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("N5").Select()
